$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new values + strip the heavy "total" style (s="10" -> none) ---
$ws.Range("B4").Value = 15.249187881000001
$ws.Range("C4").Value = 18.694921457
$ws.Range("D4").Value = 14.149005886000001
$ws.Range("E4").Value = 10.504073128000002
$ws.Range("F4").Value = 13.471785277999999
$ws.Range("G4").Value = 12.103826680000001
$ws.Range("H4").Value = 20.602144877000001
$ws.Range("I4").Value = 22.948637797
$ws.Range("J4").Value = 25.695856503000002
$ws.Range("B4:J4").Style = "Normal"

# --- Row 5 ---
$ws.Range("B5").Value = 4.2092349420000001
$ws.Range("C5").Value = 4.8672359539999999
$ws.Range("D5").Value = 2.6569098059999998
$ws.Range("E5").Value = 2.9496996360000001
$ws.Range("F5").Value = 5.0175686429999997
$ws.Range("G5").Value = 4.5489877779999999
$ws.Range("H5").Value = 6.4900402110000002
$ws.Range("I5").Value = 7.630987287
$ws.Range("J5").Value = 7.8841050680000002

# --- Row 6 ---
$ws.Range("B6").Value = 6.101510019
$ws.Range("C6").Value = 9.1498084370000008
$ws.Range("D6").Value = 6.3047681310000003
$ws.Range("E6").Value = 3.9485168939999999
$ws.Range("F6").Value = 4.6526362670000001
$ws.Range("G6").Value = 3.6893570389999999
$ws.Range("H6").Value = 8.7976447980000003
$ws.Range("I6").Value = 9.1994460629999999
$ws.Range("J6").Value = 9.2524440489999993

# --- Row 7 ---
$ws.Range("B7").Value = 1.0010529029999999
$ws.Range("C7").Value = 0.91056403900000005
$ws.Range("D7").Value = 1.2793095059999999
$ws.Range("E7").Value = 0.95995894400000004
$ws.Range("F7").Value = 1.1406763470000001
$ws.Range("G7").Value = 0.99728770099999997
$ws.Range("H7").Value = 0.73130381
$ws.Range("I7").Value = 1.1505956959999999
$ws.Range("J7").Value = 2.496598257

# --- Row 8 ---
$ws.Range("B8").Value = 3.935390017
$ws.Range("C8").Value = 3.7673130270000001
$ws.Range("D8").Value = 3.9058184429999998
$ws.Range("E8").Value = 2.6438976540000003
$ws.Range("F8").Value = 2.6609040209999999
$ws.Range("G8").Value = 2.7991441620000002
$ws.Range("H8").Value = 4.5570278310000001
$ws.Range("I8").Value = 4.9658633170000002
$ws.Range("J8").Value = 6.0627091290000008

# --- Row 12 ---
$ws.Range("B12").Value = 210.31004203399996
$ws.Range("C12").Value = 235.21579864099999
$ws.Range("D12").Value = 237.074402177
$ws.Range("E12").Value = 197.04260851699999
$ws.Range("F12").Value = 334.98417319199996
$ws.Range("G12").Value = 292.37527859700003
$ws.Range("H12").Value = 341.20708119999995
$ws.Range("I12").Value = 344.453927615
$ws.Range("J12").Value = 484.96701131999998

# --- Row 13 ---
$ws.Range("B13").Value = 191.43732836199999
$ws.Range("C13").Value = 216.479694747
$ws.Range("D13").Value = 208.19767995699999
$ws.Range("E13").Value = 167.57198866100001
$ws.Range("F13").Value = 291.77875828399999
$ws.Range("G13").Value = 260.16175313100001
$ws.Range("H13").Value = 287.49309342599997
$ws.Range("I13").Value = 278.157520347
$ws.Range("J13").Value = 413.946043062

# --- Row 14 ---
$ws.Range("B14").Value = 16.928246262000002
$ws.Range("C14").Value = 16.641993894000002
$ws.Range("D14").Value = 26.601332221999996
$ws.Range("E14").Value = 26.607125858
$ws.Range("F14").Value = 39.669487910999997
$ws.Range("G14").Value = 29.148418471000003
$ws.Range("H14").Value = 48.194645365
$ws.Range("I14").Value = 59.644049134000007
$ws.Range("J14").Value = 67.189038262000011

# --- View: move the active selection from B15:J15 up to B13:J15 ---
$ws.Range("B13:J15").Select() | Out-Null
